$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: NAQUICHE SILVA MARIA LIZETH -> 128
$ws.Range("B3").Value = 128

# Update row 4: CUBAS GARCIA ROSA ANITA -> 90
$ws.Range("B4").Value = 90

# Update row 5: MANOSALVA RUIZ SANDRA KAROLINE -> 87
$ws.Range("B5").Value = 87

# Row 6 becomes CORAS QUISPE JORGE AMERICO = 87
$ws.Range("A6").Value = "CORAS QUISPE JORGE AMERICO"
$ws.Range("B6").Value = 87

# Row 7 becomes BECERRA ASMAT CAROL STEFANY = 76
$ws.Range("A7").Value = "BECERRA ASMAT CAROL STEFANY"
$ws.Range("B7").Value = 76

# Row 8 becomes PACHECO ALISON = 69
$ws.Range("A8").Value = "PACHECO ALISON"
$ws.Range("B8").Value = 69

# Row 9 becomes CASTREJON TELLO GRECIA = 62
$ws.Range("A9").Value = "CASTREJON TELLO GRECIA"
$ws.Range("B9").Value = 62

# Row 10: SAUCEDO CABRERA CARLOS ALEXANDER -> 61
$ws.Range("B10").Value = 61
